$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header text updates: report volume/number and week-covering date range.
# ---------------------------------------------------------------------------

# "Volume 29   Number  42" -> "...43" (cell A8). Use Characters() so only the
# trailing issue-number run's text is touched.
$a8 = $ws.Range("A8")
$a8Text = $a8.Text
$oldNum = "42"
$newNum = "43"
$idx = $a8Text.LastIndexOf($oldNum)
if ($idx -ge 0) {
    $a8.Characters($idx + 1, $oldNum.Length).Text = $newNum
}

# "Report Covering the Week  10/17/2022  Through  10/23/2022" -> next week.
$c9 = $ws.Range("C9")
$c9Text = $c9.Text
$dateMap = @{
    "10/17/2022" = "10/24/2022";
    "10/23/2022" = "10/30/2022";
}
foreach ($oldDate in $dateMap.Keys) {
    $newDate = $dateMap[$oldDate]
    $curText = $c9.Text
    $dIdx = $curText.IndexOf($oldDate)
    if ($dIdx -ge 0) {
        $c9.Characters($dIdx + 1, $oldDate.Length).Text = $newDate
    }
}

# ---------------------------------------------------------------------------
# 2. Weekly crime-complaint table (rows 14-30): new week's figures.
# ---------------------------------------------------------------------------

$rowData = @(
    @{Addr='A14'; Kind='s'; Value='Murder'},
    @{Addr='C14'; Kind='n'; Value=2},
    @{Addr='D14'; Kind='n'; Value=3},
    @{Addr='E14'; Kind='n'; Value=-33.333333333333},
    @{Addr='F14'; Kind='n'; Value=7},
    @{Addr='G14'; Kind='n'; Value=7},
    @{Addr='H14'; Kind='n'; Value=0},
    @{Addr='I14'; Kind='n'; Value=66},
    @{Addr='J14'; Kind='n'; Value=83},
    @{Addr='K14'; Kind='n'; Value=-20.481927710843},
    @{Addr='L14'; Kind='n'; Value=-32.653061224489},
    @{Addr='M14'; Kind='n'; Value=-44.067796610169},
    @{Addr='N14'; Kind='n'; Value=-83.663366336633},
    @{Addr='A15'; Kind='s'; Value='Rape'},
    @{Addr='C15'; Kind='n'; Value=3},
    @{Addr='D15'; Kind='n'; Value=1},
    @{Addr='E15'; Kind='n'; Value=200},
    @{Addr='F15'; Kind='n'; Value=21},
    @{Addr='G15'; Kind='n'; Value=17},
    @{Addr='H15'; Kind='n'; Value=23.529411764705},
    @{Addr='I15'; Kind='n'; Value=216},
    @{Addr='J15'; Kind='n'; Value=181},
    @{Addr='K15'; Kind='n'; Value=19.337016574585},
    @{Addr='L15'; Kind='n'; Value=14.893617021276},
    @{Addr='M15'; Kind='n'; Value=15.508021390374},
    @{Addr='N15'; Kind='n'; Value=-57.396449704142},
    @{Addr='A16'; Kind='s'; Value='Robbery'},
    @{Addr='C16'; Kind='n'; Value=34},
    @{Addr='D16'; Kind='n'; Value=51},
    @{Addr='E16'; Kind='n'; Value=-33.333333333333},
    @{Addr='F16'; Kind='n'; Value=167},
    @{Addr='G16'; Kind='n'; Value=194},
    @{Addr='H16'; Kind='n'; Value=-13.917525773195},
    @{Addr='I16'; Kind='n'; Value=2135},
    @{Addr='J16'; Kind='n'; Value=1704},
    @{Addr='K16'; Kind='n'; Value=25.293427230046},
    @{Addr='L16'; Kind='n'; Value=19.207146845337},
    @{Addr='M16'; Kind='n'; Value=-29.584432717678},
    @{Addr='N16'; Kind='n'; Value=-84.62590912364},
    @{Addr='A17'; Kind='s'; Value='Fel. Assault'},
    @{Addr='C17'; Kind='n'; Value=71},
    @{Addr='D17'; Kind='n'; Value=57},
    @{Addr='E17'; Kind='n'; Value=24.561403508771},
    @{Addr='F17'; Kind='n'; Value=301},
    @{Addr='G17'; Kind='n'; Value=302},
    @{Addr='H17'; Kind='n'; Value=-0.331125827814},
    @{Addr='I17'; Kind='n'; Value=3453},
    @{Addr='J17'; Kind='n'; Value=2971},
    @{Addr='K17'; Kind='n'; Value=16.22349377314},
    @{Addr='L17'; Kind='n'; Value=19.937478291073},
    @{Addr='M17'; Kind='n'; Value=23.409578270193},
    @{Addr='N17'; Kind='n'; Value=-52.293451229621},
    @{Addr='A18'; Kind='s'; Value='Burglary'},
    @{Addr='C18'; Kind='n'; Value=54},
    @{Addr='D18'; Kind='n'; Value=43},
    @{Addr='E18'; Kind='n'; Value=25.581395348837},
    @{Addr='F18'; Kind='n'; Value=203},
    @{Addr='G18'; Kind='n'; Value=193},
    @{Addr='H18'; Kind='n'; Value=5.181347150259},
    @{Addr='I18'; Kind='n'; Value=1990},
    @{Addr='J18'; Kind='n'; Value=1689},
    @{Addr='K18'; Kind='n'; Value=17.821195973949},
    @{Addr='L18'; Kind='n'; Value=-11.398040961709},
    @{Addr='M18'; Kind='n'; Value=-25.412293853073},
    @{Addr='N18'; Kind='n'; Value=-80.345679012345},
    @{Addr='A19'; Kind='s'; Value='Gr. Larceny'},
    @{Addr='C19'; Kind='n'; Value=127},
    @{Addr='D19'; Kind='n'; Value=114},
    @{Addr='E19'; Kind='n'; Value=11.403508771929},
    @{Addr='F19'; Kind='n'; Value=511},
    @{Addr='G19'; Kind='n'; Value=449},
    @{Addr='H19'; Kind='n'; Value=13.80846325167},
    @{Addr='I19'; Kind='n'; Value=4973},
    @{Addr='J19'; Kind='n'; Value=3842},
    @{Addr='K19'; Kind='n'; Value=29.437792816241},
    @{Addr='L19'; Kind='n'; Value=30.490684859616},
    @{Addr='M19'; Kind='n'; Value=39.534231200897},
    @{Addr='N19'; Kind='n'; Value=-13.226313034374},
    @{Addr='A20'; Kind='s'; Value='G.L.A.'},
    @{Addr='C20'; Kind='n'; Value=42},
    @{Addr='D20'; Kind='n'; Value=31},
    @{Addr='E20'; Kind='n'; Value=35.483870967741},
    @{Addr='F20'; Kind='n'; Value=148},
    @{Addr='G20'; Kind='n'; Value=134},
    @{Addr='H20'; Kind='n'; Value=10.447761194029},
    @{Addr='I20'; Kind='n'; Value=1543},
    @{Addr='J20'; Kind='n'; Value=1275},
    @{Addr='K20'; Kind='n'; Value=21.019607843137},
    @{Addr='L20'; Kind='n'; Value=40.528233151184},
    @{Addr='M20'; Kind='n'; Value=28.90559732665},
    @{Addr='N20'; Kind='n'; Value=-80.398882113821},
    @{Addr='A21'; Kind='s'; Value='TOTAL'},
    @{Addr='C21'; Kind='n'; Value=333},
    @{Addr='D21'; Kind='n'; Value=300},
    @{Addr='E21'; Kind='n'; Value=11},
    @{Addr='F21'; Kind='n'; Value=1358},
    @{Addr='G21'; Kind='n'; Value=1296},
    @{Addr='H21'; Kind='n'; Value=4.783950617283},
    @{Addr='I21'; Kind='n'; Value=14376},
    @{Addr='J21'; Kind='n'; Value=11745},
    @{Addr='K21'; Kind='n'; Value=22.401021711366},
    @{Addr='L21'; Kind='n'; Value=18.702006440426},
    @{Addr='M21'; Kind='n'; Value=5.986434680035},
    @{Addr='N21'; Kind='n'; Value=-68.586662005069},
    @{Addr='A22'; Kind='s'; Value='Transit'},
    @{Addr='C22'; Kind='n'; Value=4},
    @{Addr='D22'; Kind='n'; Value=12},
    @{Addr='E22'; Kind='n'; Value=-66.666666666666},
    @{Addr='F22'; Kind='n'; Value=25},
    @{Addr='G22'; Kind='n'; Value=29},
    @{Addr='H22'; Kind='n'; Value=-13.793103448275},
    @{Addr='I22'; Kind='n'; Value=285},
    @{Addr='J22'; Kind='n'; Value=219},
    @{Addr='K22'; Kind='n'; Value=30.136986301369},
    @{Addr='L22'; Kind='n'; Value=4.395604395604},
    @{Addr='M22'; Kind='n'; Value=-20.833333333333},
    @{Addr='N22'; Kind='s'; Value='***.*'},
    @{Addr='A23'; Kind='s'; Value='Housing'},
    @{Addr='C23'; Kind='n'; Value=33},
    @{Addr='D23'; Kind='n'; Value=28},
    @{Addr='E23'; Kind='n'; Value=17.857142857142},
    @{Addr='F23'; Kind='n'; Value=109},
    @{Addr='G23'; Kind='n'; Value=117},
    @{Addr='H23'; Kind='n'; Value=-6.837606837606},
    @{Addr='I23'; Kind='n'; Value=1265},
    @{Addr='J23'; Kind='n'; Value=1202},
    @{Addr='K23'; Kind='n'; Value=5.241264559068},
    @{Addr='L23'; Kind='n'; Value=9.809027777777},
    @{Addr='M23'; Kind='n'; Value=28.426395939086},
    @{Addr='N23'; Kind='s'; Value='***.*'},
    @{Addr='A24'; Kind='s'; Value='Petit Larceny'},
    @{Addr='C24'; Kind='n'; Value=293},
    @{Addr='D24'; Kind='n'; Value=216},
    @{Addr='E24'; Kind='n'; Value=35.648148148148},
    @{Addr='F24'; Kind='n'; Value=1122},
    @{Addr='G24'; Kind='n'; Value=944},
    @{Addr='H24'; Kind='n'; Value=18.855932203389},
    @{Addr='I24'; Kind='n'; Value=11231},
    @{Addr='J24'; Kind='n'; Value=8525},
    @{Addr='K24'; Kind='n'; Value=31.741935483871},
    @{Addr='L24'; Kind='n'; Value=24.37430786268},
    @{Addr='M24'; Kind='n'; Value=28.914141414141},
    @{Addr='N24'; Kind='s'; Value='***.*'},
    @{Addr='A25'; Kind='s'; Value='Misd. Assault'},
    @{Addr='C25'; Kind='n'; Value=118},
    @{Addr='D25'; Kind='n'; Value=111},
    @{Addr='E25'; Kind='n'; Value=6.306306306306},
    @{Addr='F25'; Kind='n'; Value=416},
    @{Addr='G25'; Kind='n'; Value=447},
    @{Addr='H25'; Kind='n'; Value=-6.935123042505},
    @{Addr='I25'; Kind='n'; Value=4936},
    @{Addr='J25'; Kind='n'; Value=3937},
    @{Addr='K25'; Kind='n'; Value=25.374650749301},
    @{Addr='L25'; Kind='n'; Value=28.40790842872},
    @{Addr='M25'; Kind='n'; Value=-26.273338312173},
    @{Addr='N25'; Kind='s'; Value='***.*'},
    @{Addr='A26'; Kind='s'; Value='UCR Rape*'},
    @{Addr='C26'; Kind='n'; Value=8},
    @{Addr='D26'; Kind='n'; Value=8},
    @{Addr='E26'; Kind='n'; Value=0},
    @{Addr='F26'; Kind='n'; Value=34},
    @{Addr='G26'; Kind='n'; Value=31},
    @{Addr='H26'; Kind='n'; Value=9.677419354838},
    @{Addr='I26'; Kind='n'; Value=326},
    @{Addr='J26'; Kind='n'; Value=317},
    @{Addr='K26'; Kind='n'; Value=2.839116719242},
    @{Addr='L26'; Kind='n'; Value=10.884353741496},
    @{Addr='M26'; Kind='s'; Value='***.*'},
    @{Addr='N26'; Kind='s'; Value='***.*'},
    @{Addr='A27'; Kind='s'; Value='Other Sex Crimes'},
    @{Addr='C27'; Kind='n'; Value=20},
    @{Addr='D27'; Kind='n'; Value=14},
    @{Addr='E27'; Kind='n'; Value=42.857142857142},
    @{Addr='F27'; Kind='n'; Value=58},
    @{Addr='G27'; Kind='n'; Value=54},
    @{Addr='H27'; Kind='n'; Value=7.407407407407},
    @{Addr='I27'; Kind='n'; Value=522},
    @{Addr='J27'; Kind='n'; Value=571},
    @{Addr='K27'; Kind='n'; Value=-8.581436077057},
    @{Addr='L27'; Kind='n'; Value=4.609218436873},
    @{Addr='M27'; Kind='s'; Value='***.*'},
    @{Addr='N27'; Kind='s'; Value='***.*'},
    @{Addr='A28'; Kind='s'; Value='Shooting Vic.'},
    @{Addr='C28'; Kind='n'; Value=6},
    @{Addr='D28'; Kind='n'; Value=6},
    @{Addr='E28'; Kind='n'; Value=0},
    @{Addr='F28'; Kind='n'; Value=24},
    @{Addr='G28'; Kind='n'; Value=23},
    @{Addr='H28'; Kind='n'; Value=4.347826086956},
    @{Addr='I28'; Kind='n'; Value=299},
    @{Addr='J28'; Kind='n'; Value=360},
    @{Addr='K28'; Kind='n'; Value=-16.944444444444},
    @{Addr='L28'; Kind='n'; Value=-31.578947368421},
    @{Addr='M28'; Kind='n'; Value=-33.1096196868},
    @{Addr='N28'; Kind='n'; Value=-81.723716381418},
    @{Addr='A29'; Kind='s'; Value='Shooting Inc.'},
    @{Addr='C29'; Kind='n'; Value=4},
    @{Addr='D29'; Kind='n'; Value=6},
    @{Addr='E29'; Kind='n'; Value=-33.333333333333},
    @{Addr='F29'; Kind='n'; Value=22},
    @{Addr='G29'; Kind='n'; Value=22},
    @{Addr='H29'; Kind='n'; Value=0},
    @{Addr='I29'; Kind='n'; Value=250},
    @{Addr='J29'; Kind='n'; Value=291},
    @{Addr='K29'; Kind='n'; Value=-14.089347079037},
    @{Addr='L29'; Kind='n'; Value=-30.939226519337},
    @{Addr='M29'; Kind='n'; Value=-31.129476584022},
    @{Addr='N29'; Kind='n'; Value=-82.970027247956},
    @{Addr='A30'; Kind='s'; Value='Hate Crimes'},
    @{Addr='C30'; Kind='s'; Value='***.*'},
    @{Addr='D30'; Kind='n'; Value=3},
    @{Addr='E30'; Kind='n'; Value=-100},
    @{Addr='F30'; Kind='n'; Value=4},
    @{Addr='G30'; Kind='n'; Value=5},
    @{Addr='H30'; Kind='n'; Value=-20},
    @{Addr='I30'; Kind='n'; Value=71},
    @{Addr='J30'; Kind='n'; Value=54},
    @{Addr='K30'; Kind='n'; Value=31.481481481481},
    @{Addr='L30'; Kind='n'; Value=108.823529411765},
    @{Addr='M30'; Kind='s'; Value='***.*'},
    @{Addr='N30'; Kind='s'; Value='***.*'},
)


foreach ($item in $rowData) {
    $cell = $ws.Range($item.Addr)
    if ($item.Kind -eq 's') {
        $cell.Value2 = $item.Value
    } else {
        $cell.Value2 = [double]$item.Value
    }
}

# ---------------------------------------------------------------------------
# 3. Style fix-ups where a cell's data type changed (text <-> number) so the
#    new value keeps the same look as its column peers.
# ---------------------------------------------------------------------------

$countFormat = "#,##0"
$pctFormat = '#,##0.0;"-"#,##0.0'

$ws.Range("C14").NumberFormat = $countFormat
$ws.Range("D14").NumberFormat = $countFormat
$ws.Range("E14").NumberFormat = $pctFormat

$ws.Range("C30").NumberFormat = "General"
$ws.Range("C30").HorizontalAlignment = $ws.Range("D30").HorizontalAlignment
$ws.Range("C30").VerticalAlignment = $ws.Range("D30").VerticalAlignment
